$d = $word.ActiveDocument

# 1. Merge the "verbal" / "and written communication skills" runs back into one
#    continuous sentence (removing the split that held the _GoBack bookmark).
$d.Content.Find.Execute(
    "high level of accuracy and strong attention to detail, with excellent verbal and written communication skills.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "high level of accuracy and strong attention to detail, with excellent verbal and written communication skills.",
    2) | Out-Null

# 2. "Python" -> "C,  C++" in the Proficient computer-competencies line.
$d.Content.Find.Execute(
    "Proficient: Javascript, HTML, MySQL, node.js, PHP, Ruby, Scala, Python, and jQuery",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Proficient: Javascript, HTML, MySQL, node.js, PHP, Ruby, Scala, C,  C++, and jQuery",
    2) | Out-Null

# 3. "Familiar: C#, C++" -> "Familiar: C#, Python"
$d.Content.Find.Execute(
    "Familiar: C#, C++",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Familiar: C#, Python",
    2) | Out-Null

# 4. "EMPLOYMENT HISTORY" heading re-typed (splits into two runs with the
#    _GoBack bookmark now sitting between "EMPLOYMENT " and "HISTORY").
$d.Content.Find.Execute(
    "EMPLOYMENT HISTORY",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "EMPLOYMENT HISTORY",
    2) | Out-Null
